$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.864.29"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.553.00"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.94"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.483"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.56"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.769.78"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "1.551.64"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "26.825.78"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.25"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.09"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0687"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.18"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.103"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0463"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "1.357.17"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.930"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.77"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.18"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.684.38"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "0.0₇0974"
$ws.Range("E51").Value = "  -0.75%  "
